# "add header to menus"
# The GUI sheet tracks, per menu class (column A), the assigned teammate
# (column B) and the fraction of work credited towards a "header" task
# (column C). A "header" item was added for every menu that Soheil owns,
# so his rows' progress (column C) now reflect that extra completed task.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GUI")

# LoginMenu and MainMenu are now fully credited (header was their only task).
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 1

# The remaining Soheil-owned menus now have the header task counted as half
# of their work.
$ws.Range("C6").Value = 0.5
$ws.Range("C9").Value = 0.5
$ws.Range("C13").Value = 0.5
$ws.Range("C20").Value = 0.5
$ws.Range("C23").Value = 0.5
$ws.Range("C24").Value = 0.5
$ws.Range("C27").Value = 0.5

# Leave the selection where the last edit was made.
[void]$ws.Range("C27").Select()
